# daily auto push: 2026-02-25 05:03 UTC
# Insert one new log row for 2026/02/25 (hour 13) just after the existing
# 2026/02/25 entries (rows 868-870) and before the 2026/12/29 block, pushing
# every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("871:871").Insert()

# Write the date as literal text (leading apostrophe forces text so Excel
# doesn't reinterpret "2026/02/25" as a date serial), matching how the rest
# of column A is stored, then clear the resulting quote-prefix style so the
# cell carries no extra formatting.
$ws.Range("A871").Value = "'2026/02/25"
$ws.Range("A871").Style = "Normal"

$ws.Range("B871").Value = "水"
$ws.Range("C871").Value = 13
$ws.Range("D871").Value = 201
